$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.154.40"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.912.24"
$ws.Range("E3").Value = "  -3.49%  "
$ws.Range("E4").Value = "  -1.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.93"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("E7").Value = "  -6.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3999"
$ws.Range("E8").Value = "  -4.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.27"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08411"
$ws.Range("E10").Value = "  -9.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.043"
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.901.05"
$ws.Range("E13").Value = "  -4.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.418"
$ws.Range("E14").Value = "  -6.56%  "
$ws.Range("E15").Value = "  -6.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.50"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001058"
$ws.Range("E18").Value = "  -4.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06577"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.96"
$ws.Range("E20").Value = "  -6.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.709"
$ws.Range("E22").Value = "  -4.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.145.46"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  -5.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.277"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.130.05"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.16"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -4.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.126"
$ws.Range("E29").Value = "  -6.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.698"
$ws.Range("E30").Value = "  -9.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.94"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9725"
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09569"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.442"
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.629"
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.531"
$ws.Range("E36").Value = "  -5.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.815"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02292"
$ws.Range("E38").Value = "  -5.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06150"
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.229"
$ws.Range("E40").Value = "  -6.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6123"
$ws.Range("E41").Value = "  -5.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.98"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1900"
$ws.Range("E44").Value = "  -5.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.302"
$ws.Range("E45").Value = "  -4.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5841"
$ws.Range("E46").Value = "  -6.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.70"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.021"
$ws.Range("E48").Value = "  -7.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.430"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06850"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.46"
$ws.Range("E51").Value = "  -3.49%  "
